# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2  = 2
    3  = 1
    4  = 0
    7  = 0
    8  = 1
    9  = 0
    10 = 2
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 2
    17 = 3
    18 = 2
    19 = 1
    20 = 1
    21 = 2
    22 = 2
    23 = 2
    24 = 0
    25 = 3
    26 = 0
    27 = 0
    28 = 0
    29 = 2
    30 = 4
    31 = 2
    32 = 0
    33 = 0
    34 = 2
    35 = 1
    36 = 0
    38 = 1
    39 = 2
    40 = 2
    42 = 2
    43 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
